# Replace the final "bare URL" paragraph (seamless-metal-door-texture) with
# the full "Barricade Texture / Grass Texture / Skybox" asset listing blocks,
# plus the start of a new "vehicles-assets-pt1" entry, per the commit:
# "barricades fixed, skybox, car mesh started".

$d = $word.ActiveDocument

# The paragraph that currently just contains the bare door-texture URL is the
# last paragraph in the document; replace its content (and append all of the
# new paragraphs that follow it) via a single WordOpenXML fragment. Using
# InsertXML keeps full control of run/paragraph-mark formatting (bold +
# bCs on both the heading run and its paragraph mark) instead of relying on
# Font property side effects.
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range

$xmlFragment = @"
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Barricade Texture</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>URL:</w:t></w:r>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:r><w:t>https://opengameart.org/content/seamless-metal-door-texture</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>Date of Download: 01/04/2023</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>License: Creative Commons</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
  <w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Grass Texture:</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t xml:space="preserve">URL: </w:t></w:r>
  <w:r><w:t>https://opengameart.org/content/grass-texture-0</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t xml:space="preserve">Date Of Download: </w:t></w:r>
  <w:r><w:t>03/04/2023</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>License: Creative Commons</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:t>Skybox:</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t xml:space="preserve">URL: </w:t></w:r>
  <w:r><w:t>https://opengameart.org/content/mountain-skyboxes</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>Date of Download: 03/04/2023</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:tab/><w:t>License: Creative Commons</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>
  <w:r><w:t>https://opengameart.org/content/vehicles-assets-pt1</w:t></w:r>
</w:p>
<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'/>
"@

$null = $insertRange.InsertXML($xmlFragment)

# Turn the three plain-text URLs that were just inserted into real
# hyperlinks (rStyle="Hyperlink"), mirroring how the existing Tunnel/Track
# Texture entries in the document already work.
$urls = @(
  "https://opengameart.org/content/seamless-metal-door-texture",
  "https://opengameart.org/content/grass-texture-0",
  "https://opengameart.org/content/mountain-skyboxes"
)

foreach ($u in $urls) {
  $searchRange = $d.Content
  $found = $searchRange.Find.Execute($u, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
  if ($found) {
    $null = $d.Hyperlinks.Add($searchRange, $u)
  }
}
